{"js": "// Repair corruptions: remove the stray w:pStyle (\"Header\"/\"Footer\") that\n// was left in every header/footer paragraph's pPr. All header/footer\n// parts in this document consist of a single empty paragraph whose\n// paragraph properties only carry that pStyle, so clearing the\n// paragraph style collapses the (now empty) pPr away entirely -- which\n// is the same net effect as deleting just the <w:pStyle/> child.\n\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst headerFooterTypes = [\"Primary\", \"FirstPage\", \"EvenPages\"];\n\nconst bodies = [];\nfor (let i = 0; i < sections.items.length; i++) {\n  const sec = sections.items[i];\n  for (const type of headerFooterTypes) {\n    bodies.push(sec.getHeader(type));\n    bodies.push(sec.getFooter(type));\n  }\n}\n\nfor (const body of bodies) {\n  body.paragraphs.load(\"items\");\n}\nawait context.sync();\n\nfor (const body of bodies) {\n  for (const p of body.paragraphs.items) {\n    // Resetting to the built-in \"Normal\" style removes the w:pStyle\n    // element (and, since nothing else lives in pPr here, the pPr\n    // itself disappears too -- functionally identical to just\n    // stripping the pStyle child as the diff does).\n    p.style = \"Normal\";\n  }\n}\nawait context.sync();\n", "ps1": "# Repair corruptions: remove the stray pStyle (\"Header\"/\"Footer\") that\n# was left in every header/footer paragraph. Every header/footer story\n# in this document is a single empty paragraph whose only formatting is\n# that paragraph style, so resetting the paragraph style to the builtin\n# \"Normal\" style clears the w:pStyle element (and, since nothing else\n# lives in its pPr, the now-empty pPr collapses away too) -- the same\n# net effect as deleting just the <w:pStyle/> child.\n\n$d = $word.ActiveDocument\n\nfor ($i = 1; $i -le $d.Sections.Count; $i++) {\n    $sec = $d.Sections.Item($i)\n\n    for ($j = 1; $j -le $sec.Headers.Count; $j++) {\n        $hdr = $sec.Headers.Item($j)\n        # Apply the style reset directly on the header's Range (rather\n        # than hopping through .Range.Paragraphs) so the edit stays\n        # anchored to this header/footer story instead of the body.\n        $hdr.Range.Style = \"Normal\"\n    }\n\n    for ($j = 1; $j -le $sec.Footers.Count; $j++) {\n        $ftr = $sec.Footers.Item($j)\n        $ftr.Range.Style = \"Normal\"\n    }\n}\n"}
